$wb = $excel.ActiveWorkbook

# ---- Sheet: Short Term ----
$ws1 = $wb.Worksheets.Item("Short Term")

# Update existing rows 119-124 (columns B:G)
$shortTermUpdates = @{
    119 = @(-8.38, -5.61, -2.69, 21.73, 20.81, 1.06)
    120 = @(5.11, -2.61, 7.35, 26.86, 31.34, -16.15)
    121 = @(-8.99, -10.02, -3.08, 15.27, 17.27, -1.23)
    122 = @(37.07, 37.42, -0.49, 39.65, 36.83, 10.82)
    123 = @(0.05, 4.08, 21.28, 66.87, 58.24, 22.15)
    124 = @(36.75, 38.11, -7.73, 104.88, 98.68, 11.25)
}

foreach ($row in $shortTermUpdates.Keys) {
    $vals = $shortTermUpdates[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i  # B = 2
        $ws1.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# Add new row 125
$ws1.Cells.Item(124, 1).Copy()
$ws1.Cells.Item(125, 1).PasteSpecial(-4122)
$ws1.Cells.Item(125, 1).Value = 45748
$ws1.Cells.Item(125, 2).Value = -36.96
$ws1.Cells.Item(125, 3).Value = -39.3
$ws1.Cells.Item(125, 4).Value = -11.95
$ws1.Cells.Item(125, 5).Value = 17.77
$ws1.Cells.Item(125, 6).Value = 12.65
$ws1.Cells.Item(125, 7).Value = -0.17

# ---- Sheet: Medium Term ----
$ws2 = $wb.Worksheets.Item("Medium Term")

$mediumTermUpdates = @{
    105 = @(22.71, 13.55, 9.2)
    106 = @(30.38, 17.23, 12.33)
    107 = @(23.35, 22.3, 14.43)
    108 = @(26.08, 15.96)
    109 = @(34.31, 21.09)
    110 = @(65.1, 45.6, 30.33)
}

foreach ($row in $mediumTermUpdates.Keys) {
    $vals = $mediumTermUpdates[$row]
    if ($row -eq 108 -or $row -eq 109) {
        # Only columns C and D change for these rows
        $ws2.Cells.Item($row, 3).Value = $vals[0]
        $ws2.Cells.Item($row, 4).Value = $vals[1]
    } else {
        for ($i = 0; $i -lt $vals.Length; $i++) {
            $col = 2 + $i  # B = 2
            $ws2.Cells.Item($row, $col).Value = $vals[$i]
        }
    }
}

# Add new row 111
$ws2.Cells.Item(110, 1).Copy()
$ws2.Cells.Item(111, 1).PasteSpecial(-4122)
$ws2.Cells.Item(111, 1).Value = 45748
$ws2.Cells.Item(111, 2).Value = 56.23
$ws2.Cells.Item(111, 3).Value = 43.61
$ws2.Cells.Item(111, 4).Value = 28.81
